# Scheduled-runner profit recalculation: refresh market-board derived
# columns (H..N: currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit
# NQ/HQ) for the leve rows whose source prices moved since the last sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4233.3335
$ws.Range("I74").Value = 4280
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 4280
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -3344
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 4233.3335
$ws.Range("I77").Value = 4280
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 21400
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -16720
$ws.Range("N77").Value = -29360
$ws.Range("H135").Value = 4649
$ws.Range("I135").Value = 4649
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 41841
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -39306
$ws.Range("H138").Value = 9849.447
$ws.Range("I138").Value = 7425.1816
$ws.Range("J138").Value = 10837.111
$ws.Range("K138").Value = 22275.5448
$ws.Range("L138").Value = 32511.333
$ws.Range("M138").Value = -17135.5448
$ws.Range("N138").Value = -42791.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5946.2856
$ws.Range("I32").Value = 4875.758
$ws.Range("J32").Value = 23610
$ws.Range("K32").Value = 4875.758
$ws.Range("L32").Value = 23610
$ws.Range("M32").Value = -4588.758
$ws.Range("N32").Value = -24184
$ws.Range("H61").Value = 5770.857
$ws.Range("I61").Value = 1699
$ws.Range("J61").Value = 7399.6
$ws.Range("K61").Value = 1699
$ws.Range("L61").Value = 7399.6
$ws.Range("M61").Value = -1487
$ws.Range("N61").Value = -7823.6
$ws.Range("H74").Value = 2270.5
$ws.Range("I74").Value = 2270.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2270.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1396.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2270.5
$ws.Range("I77").Value = 2270.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 11352.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -6984.5
$ws.Range("N77").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 1243.8667
$ws.Range("I132").Value = 1272.3572
$ws.Range("J132").Value = 845
$ws.Range("K132").Value = 3817.0716
$ws.Range("L132").Value = 2535
$ws.Range("M132").Value = -1287.0716
$ws.Range("N132").Value = -7595
$ws.Range("H136").Value = 5770.857
$ws.Range("I136").Value = 1699
$ws.Range("J136").Value = 7399.6
$ws.Range("K136").Value = 5097
$ws.Range("L136").Value = 22198.8
$ws.Range("M136").Value = -2547
$ws.Range("N136").Value = -27298.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2215.4375
$ws.Range("I94").Value = 1419.0769
$ws.Range("J94").Value = 5666.3335
$ws.Range("K94").Value = 1419.0769
$ws.Range("L94").Value = 5666.3335
$ws.Range("M94").Value = -968.0769
$ws.Range("N94").Value = -6568.3335
$ws.Range("H99").Value = 1999.75
$ws.Range("I99").Value = 1999.6666
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1999.6666
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -501.6666
$ws.Range("N99").Value = -4996
$ws.Range("H106").Value = 38000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 38000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 38000
$ws.Range("N106").Value = -40524
$ws.Range("H134").Value = 1852.25
$ws.Range("I134").Value = 1852.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5556.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3021.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 457.75
$ws.Range("I6").Value = 443.66666
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 443.66666
$ws.Range("L6").Value = 500
$ws.Range("M6").Value = -330.66666
$ws.Range("N6").Value = -726
$ws.Range("H7").Value = 1917.5714
$ws.Range("I7").Value = 2121.4
$ws.Range("J7").Value = 1804.3334
$ws.Range("K7").Value = 2121.4
$ws.Range("L7").Value = 1804.3334
$ws.Range("M7").Value = -2008.4
$ws.Range("N7").Value = -2030.3334
$ws.Range("H17").Value = 352.5
$ws.Range("I17").Value = 205
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 205
$ws.Range("L17").Value = 500
$ws.Range("M17").Value = -31
$ws.Range("N17").Value = -848
$ws.Range("H25").Value = 1733.3334
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 1850
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1850
$ws.Range("M25").Value = -1326
$ws.Range("N25").Value = -2198
$ws.Range("H41").Value = 3333
$ws.Range("I41").Value = 3333
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3333
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2905
$ws.Range("N41").ClearContents()
$ws.Range("H140").Value = 500000
$ws.Range("I140").Value = 500000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 500000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -494820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44
$ws.Range("I2").Value = 44
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 264
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -151
$ws.Range("N2").ClearContents()
$ws.Range("H46").Value = 156.66667
$ws.Range("I46").Value = 185.25
$ws.Range("J46").Value = 99.5
$ws.Range("K46").Value = 555.75
$ws.Range("L46").Value = 298.5
$ws.Range("M46").Value = -464.75
$ws.Range("N46").Value = -480.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 25000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576
$ws.Range("H80").Value = 3258.6
$ws.Range("I80").Value = 3165.6667
$ws.Range("J80").Value = 3398
$ws.Range("K80").Value = 3165.6667
$ws.Range("L80").Value = 3398
$ws.Range("M80").Value = -2167.6667
$ws.Range("N80").Value = -5394
$ws.Range("H81").Value = 25000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 25000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996
$ws.Range("H83").Value = 3258.6
$ws.Range("I83").Value = 3165.6667
$ws.Range("J83").Value = 3398
$ws.Range("K83").Value = 15828.3335
$ws.Range("L83").Value = 16990
$ws.Range("M83").Value = -10836.3335
$ws.Range("N83").Value = -26974
$ws.Range("H84").Value = 25000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 25000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2836.1875
$ws.Range("I22").Value = 775
$ws.Range("J22").Value = 3130.6428
$ws.Range("K22").Value = 775
$ws.Range("L22").Value = 3130.6428
$ws.Range("M22").Value = -480
$ws.Range("N22").Value = -3720.6428
$ws.Range("H27").Value = 2836.1875
$ws.Range("I27").Value = 775
$ws.Range("J27").Value = 3130.6428
$ws.Range("K27").Value = 775
$ws.Range("L27").Value = 3130.6428
$ws.Range("M27").Value = -668
$ws.Range("N27").Value = -3344.6428
$ws.Range("H46").Value = 1735.375
$ws.Range("I46").Value = 2450
$ws.Range("J46").Value = 1497.1666
$ws.Range("K46").Value = 2450
$ws.Range("L46").Value = 1497.1666
$ws.Range("M46").Value = -2262
$ws.Range("N46").Value = -1873.1666
$ws.Range("H57").Value = 25050000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 25050000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 25050000
$ws.Range("N57").Value = -25051132
$ws.Range("H61").Value = 2861.4285
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 3005
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 3005
$ws.Range("M61").Value = -1798
$ws.Range("N61").Value = -3409
$ws.Range("H93").Value = 531.3333
$ws.Range("I93").Value = 597
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 597
$ws.Range("L93").Value = 400
$ws.Range("M93").Value = 651
$ws.Range("N93").Value = -2896
$ws.Range("H113").Value = 2861.4285
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3005
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3005
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 33749.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 33749.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 33749.5
$ws.Range("N82").Value = -34515.5
$ws.Range("H85").Value = 33749.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 33749.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 33749.5
$ws.Range("N85").Value = -36401.5
$ws.Range("H136").Value = 27460.4
$ws.Range("I136").Value = 28379.37
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 85138.11
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -82588.11
$ws.Range("N136").Value = -35100
